$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/employee-performance-rating"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# --- Elements sheet updates ---
$elements = $wb.Worksheets.Item("Elements")
# Fixed Value of Extension.url mirrors the StructureDefinition's canonical URL
$elements.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/employee-performance-rating"
# Constraint(s) column for the top-level Extension row is cleared (the
# ele-1/ext-1 constraint now only applies to the Extension.extension row)
$elements.Range("AI2").Value = ""
